$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a plain "Text" number format on the Price cells whose new value would
# otherwise be auto-parsed as a real number (losing the literal digit grouping
# / trailing-zero formatting used by this sheet), then write the literal string.
$textPriceRows = @(5,6,8,12,16,18,19,20,21,22,25,29,31,33,34,35,36,37,39,43,44,45,46,47,48,50,51)
foreach ($r in $textPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '60.460.55'
$ws.Cells.Item(2, 5).Value = '  -1.70%  '
$ws.Cells.Item(3, 4).Value = '2.898.10'
$ws.Cells.Item(3, 5).Value = '  -2.39%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Value = '525.44'
$ws.Cells.Item(5, 5).Value = '  -2.41%  '
$ws.Cells.Item(6, 4).Value = '141.75'
$ws.Cells.Item(6, 5).Value = '  -5.69%  '
$ws.Cells.Item(8, 4).Value = '0.549'
$ws.Cells.Item(8, 5).Value = '  -3.11%  '
$ws.Cells.Item(9, 4).Value = '2.903.37'
$ws.Cells.Item(9, 5).Value = '  -2.45%  '
$ws.Cells.Item(10, 5).Value = '  -4.84%  '
$ws.Cells.Item(11, 5).Value = '  -2.56%  '
$ws.Cells.Item(12, 4).Value = '0.359'
$ws.Cells.Item(12, 5).Value = '  -2.24%  '
$ws.Cells.Item(13, 4).Value = '3.407.79'
$ws.Cells.Item(13, 5).Value = '  -2.36%  '
$ws.Cells.Item(14, 5).Value = '  +2.20%  '
$ws.Cells.Item(15, 4).Value = '60.501.30'
$ws.Cells.Item(15, 5).Value = '  -1.76%  '
$ws.Cells.Item(16, 4).Value = '22.67'
$ws.Cells.Item(16, 5).Value = '  -3.89%  '
$ws.Cells.Item(17, 4).Value = '2.902.34'
$ws.Cells.Item(17, 5).Value = '  -2.50%  '
$ws.Cells.Item(18, 4).Value = '0.0000141'
$ws.Cells.Item(18, 5).Value = '  -3.95%  '
$ws.Cells.Item(19, 4).Value = '4.98'
$ws.Cells.Item(19, 5).Value = '  -3.33%  '
$ws.Cells.Item(20, 4).Value = '11.63'
$ws.Cells.Item(20, 5).Value = '  -3.12%  '
$ws.Cells.Item(21, 4).Value = '353.17'
$ws.Cells.Item(21, 5).Value = '  -7.10%  '
$ws.Cells.Item(22, 4).Value = '6.55'
$ws.Cells.Item(22, 5).Value = '  -1.60%  '
$ws.Cells.Item(23, 5).Value = '  +0.02%  '
$ws.Cells.Item(24, 5).Value = '  +0.95%  '
$ws.Cells.Item(25, 4).Value = '64.51'
$ws.Cells.Item(25, 5).Value = '  -1.38%  '
$ws.Cells.Item(26, 5).Value = '  -3.38%  '
$ws.Cells.Item(27, 5).Value = '  -5.31%  '
$ws.Cells.Item(28, 5).Value = '  +0.15%  '
$ws.Cells.Item(29, 4).Value = '7.83'
$ws.Cells.Item(29, 5).Value = '  -4.72%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0839'
$ws.Cells.Item(30, 5).Value = '  -10.62%  '
$ws.Cells.Item(31, 4).Value = '0.999'
$ws.Cells.Item(31, 5).Value = '  -0.02%  '
$ws.Cells.Item(32, 5).Value = '  -2.22%  '
$ws.Cells.Item(33, 4).Value = '19.59'
$ws.Cells.Item(33, 5).Value = '  -4.08%  '
$ws.Cells.Item(34, 4).Value = '149.75'
$ws.Cells.Item(34, 5).Value = '  -6.88%  '
$ws.Cells.Item(35, 4).Value = '4.33'
$ws.Cells.Item(35, 5).Value = '  -6.79%  '
$ws.Cells.Item(36, 4).Value = '5.58'
$ws.Cells.Item(36, 5).Value = '  -5.45%  '
$ws.Cells.Item(37, 4).Value = '0.998'
$ws.Cells.Item(37, 5).Value = '  -6.68%  '
$ws.Cells.Item(38, 5).Value = '  -4.80%  '
$ws.Cells.Item(39, 4).Value = '37.73'
$ws.Cells.Item(39, 5).Value = '  +0.64%  '
$ws.Cells.Item(40, 5).Value = '  -4.61%  '
$ws.Cells.Item(41, 5).Value = '  -4.82%  '
$ws.Cells.Item(42, 4).Value = '2.288.25'
$ws.Cells.Item(42, 5).Value = '  -4.82%  '
$ws.Cells.Item(43, 4).Value = '0.647'
$ws.Cells.Item(43, 5).Value = '  -3.13%  '
$ws.Cells.Item(44, 4).Value = '0.0581'
$ws.Cells.Item(44, 5).Value = '  -1.39%  '
$ws.Cells.Item(45, 4).Value = '20.34'
$ws.Cells.Item(45, 5).Value = '  -7.88%  '
$ws.Cells.Item(46, 4).Value = '0.998'
$ws.Cells.Item(46, 5).Value = '  +0.05%  '
$ws.Cells.Item(47, 4).Value = '4.93'
$ws.Cells.Item(47, 5).Value = '  -3.24%  '
$ws.Cells.Item(48, 4).Value = '0.0237'
$ws.Cells.Item(48, 5).Value = '  -3.76%  '
$ws.Cells.Item(49, 5).Value = '  -1.27%  '
$ws.Cells.Item(50, 4).Value = '0.0919'
$ws.Cells.Item(50, 5).Value = '  -3.20%  '
$ws.Cells.Item(51, 4).Value = '248.20'
$ws.Cells.Item(51, 5).Value = '  -6.79%  '
